# Update the "datetimeFigureOut" Date Placeholder field text from
# 08/12/2022 -> 11/12/2022 across the Slide Master and every Slide
# Layout (the diff touches this literal string in 12 places: the
# master + all 11 layouts).

$p = $ppt.ActivePresentation

$oldDate = "08/12/2022"
$newDate = "11/12/2022"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout (CustomLayouts) hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
